# Rebuild the "Year"-pivoted table (1 header row x 13 topic columns) into a
# tall 3-column table (GRI disclosure id / title / page), 14 rows total.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start from a clean sheet so no stale values/formatting from the old
# 13-column layout survive (old columns D:M must disappear entirely).
$ws.Cells.Clear()

# --- Header row (row 1) -----------------------------------------------
# Same look as the old header cells: bold font, thin box border, centered
# horizontally, top-aligned vertically.
$header = $ws.Range("A1:C1")
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108   # xlCenter
$header.VerticalAlignment = -4160     # xlTop
$header.Borders.LineStyle = 1

# The header row's contents are numeric-looking ("4.1", "30") but must be
# stored as text, like the rest of the header row ("GRI 401: Employment
# 2016"). A leading apostrophe forces text storage without touching the
# cell's number format.
$ws.Range("A1").Formula = "'4.1"
$ws.Range("B1").Value = "GRI 401: Employment 2016"
$ws.Range("C1").Formula = "'30"

# --- Data rows (rows 2-14) ---------------------------------------------
$rows = @(
    @(4.2,  "GRI 402: Labor/management relations 2016",          33),
    @(4.3,  "GRI 403: Occupational health and safety 2016",      33),
    @(4.4,  "GRI 404: Training and education 2016",               34),
    @(4.5,  "GRI 405: Diversity and equal opportunity 2016",      36),
    @(4.6,  "GRI 406: Non-discrimination 2016",                   38),
    @(4.7,  "GRI 412: Human Rights Assessment 2016",              39),
    @(4.8,  "GRI 413: Local communities 2016",                    40),
    @(4.9,  "GRI 415: Public policy 2016",                        41),
    @(4.1,  "GRI 417: Marketing and labeling 2016",               42),
    @(4.11, "GRI 418: Customer privacy 2016",                     42),
    @(4.12, "GRI 419: Socioeconomic compliance 2016",             42),
    @(5,    "Financial services sector disclosures",              43),
    @(5.1,  "Product portfolio",                                  43)
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}
